# Update the row-2 record in the Decision Summary sheet.
# Numeric-looking values (E2:I2) are forced to remain TEXT (matching the
# original inlineStr/string cell type) by using a leading apostrophe,
# the standard Excel "store as text" convention - same as setting
# NumberFormat "@" / quotePrefix, but without touching number formats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ATMA"
$ws.Range("B2").Value = "Camst"
$ws.Range("C2").Value = "Rproj_5"
$ws.Range("E2").Value = "'23"
$ws.Range("F2").Value = "'234"
$ws.Range("G2").Value = "'546"
$ws.Range("H2").Value = "'234"
$ws.Range("I2").Value = "'645"
